# "Generate Report for Handoff"
# Adds a new localization-status row (file
# 71be0cc9-dd1c-4666-a97c-2479f001881e...md, status "Ready for handoff")
# to the Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

$newFile = "71be0cc9-dd1c-4666-a97c-2479f001881eooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$newPath = "e2e\" + $newFile
$newUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f673cec4aa1df65609f724515dbea7031d92d505/e2e/" + $newFile

$zhXlf = "71be0cc9-dd1c-4666-a97c-2479f001881eoooooooooooooooooooooooooooooooooooooooo.90b8fe820a54f9bfb1538b327695d44489de658e.zh-cn.xlf"
$deXlf = "71be0cc9-dd1c-4666-a97c-2479f001881eoooooooooooooooooooooooooooooooooooooooo.90b8fe820a54f9bfb1538b327695d44489de658e.de-de.xlf"

# ---------------------------------------------------------------
# Overview sheet -> row 3
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $newFile
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newUrl, "", "", $newPath)
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-26 16:28:38"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$loOverview = $wsOverview.ListObjects.Item("Overview")
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------
# zh-cn sheet -> row 3
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $newUrl, "", "", $newFile)
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "False"
$wsZh.Range("G3").Value = $zhXlf
$wsZh.Range("H3").Value = "2016-08-26 16:28:34"
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("M3").Value = "True"
$wsZh.Range("O3").Value = "False"

$loZh = $wsZh.ListObjects.Item("zh-cn")
$loZh.Resize($wsZh.Range("A1:P3"))

# ---------------------------------------------------------------
# de-de sheet -> row 3
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $newUrl, "", "", $newFile)
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "False"
$wsDe.Range("G3").Value = $deXlf
$wsDe.Range("H3").Value = "2016-08-26 16:28:38"
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("M3").Value = "True"
$wsDe.Range("O3").Value = "False"

$loDe = $wsDe.ListObjects.Item("de-de")
$loDe.Resize($wsDe.Range("A1:P3"))

# ---------------------------------------------------------------
# Column width tweaks (widened to fit the new longer values)
# ---------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797
$wsZh.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsDe.Columns.Item(3).ColumnWidth = 17.2159881591797

Write-Host "Report generated for handoff."
